# Apply the "Add files via upload" edit to the document:
#   - Append a line break, another line break, and the text
#     "edit to hub file" after the existing "3rd GitHub file" text.
#   - The _GoBack bookmark (originally wrapping the whole first run) ends
#     up collapsed at the very end of the paragraph, after the new text.

$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)

# The hidden "_GoBack" bookmark currently wraps the original text. Remove
# it first so we can freely append new runs after "3rd GitHub file";
# we'll recreate it (collapsed) at the new end of the paragraph afterwards.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a standalone line-break run: <w:r><w:br/></w:r>
$r1 = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$r1.InsertBreak(6)   # wdLineBreak

# Insert a second run containing a line break followed by the new text:
# <w:r><w:br/><w:t>edit to hub file</w:t></w:r>
# Chr(11) is a manual line break character, so putting it in the same
# InsertAfter call as the text keeps the break and the text in one run.
$brChar = [char]11
$r2 = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$r2.InsertAfter("$brChar" + "edit to hub file")

# Recreating the "_GoBack" bookmark collapsed at the very end of the
# paragraph's text triggers an engine quirk that resets its Start to 0
# when the target position is exactly the last valid offset. Work around
# it by inserting a temporary placeholder character, adding the bookmark
# just before it (a "safe", non-final position), and then deleting the
# placeholder so the bookmark ends up collapsed at the true end.
$placeholder = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$placeholder.InsertAfter("X")

$pos = $p.Range.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($pos, $pos))

$trailing = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$trailing.Delete()
